$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 2.43753139
$ws.Cells.Item(2, 8).Value = 166.85076695
$ws.Cells.Item(2, 13).Value = 0.2808116673432284
$ws.Cells.Item(2, 14).Value = 36.96750933532022

$ws.Cells.Item(3, 7).Value = 1.43362922
$ws.Cells.Item(3, 8).Value = 159.68834131
$ws.Cells.Item(3, 13).Value = 0.3857751479011183
$ws.Cells.Item(3, 14).Value = 50.99682913882356

$ws.Cells.Item(4, 7).Value = 1.26306222
$ws.Cells.Item(4, 8).Value = 56.55318427
$ws.Cells.Item(4, 13).Value = 0.2686290274461385
$ws.Cells.Item(4, 14).Value = 21.46195678839939

$ws.Cells.Item(5, 7).Value = 0.730662
$ws.Cells.Item(5, 8).Value = 47.42466252
$ws.Cells.Item(5, 13).Value = 0.1994292735205565
$ws.Cells.Item(5, 14).Value = 18.003032351658

$ws.Cells.Item(6, 7).Value = 0.57986558
$ws.Cells.Item(6, 8).Value = 16.04977179
$ws.Cells.Item(6, 13).Value = 0.1691694069141915
$ws.Cells.Item(6, 14).Value = 7.855137124474721

$ws.Cells.Item(7, 7).Value = 0.34253878
$ws.Cells.Item(7, 8).Value = 11.79694189
$ws.Cells.Item(7, 13).Value = 0.09894841761271202
$ws.Cells.Item(7, 14).Value = 4.939485261636971

$ws.Cells.Item(8, 7).Value = 0.33216157
$ws.Cells.Item(8, 8).Value = 6.24381025
$ws.Cells.Item(8, 13).Value = 0.1163972933660317
$ws.Cells.Item(8, 14).Value = 3.833098956267941

$ws.Cells.Item(9, 7).Value = 0.19160893
$ws.Cells.Item(9, 8).Value = 4.31384833
$ws.Cells.Item(9, 13).Value = 0.05772261005937939
$ws.Cells.Item(9, 14).Value = 1.920355393550123

$ws.Cells.Item(10, 7).Value = 0.20314113
$ws.Cells.Item(10, 8).Value = 2.75909081
$ws.Cells.Item(10, 13).Value = 0.08173785510263822
$ws.Cells.Item(10, 14).Value = 1.904877334205648

$ws.Cells.Item(11, 7).Value = 0.12062729
$ws.Cells.Item(11, 8).Value = 2.02756292
$ws.Cells.Item(11, 13).Value = 0.04237322676056193
$ws.Cells.Item(11, 14).Value = 1.04764165607101

$ws.Cells.Item(12, 7).Value = 0.12852231
$ws.Cells.Item(12, 8).Value = 1.31716232
$ws.Cells.Item(12, 13).Value = 0.0569629458522298
$ws.Cells.Item(12, 14).Value = 0.9971903187449306

$ws.Cells.Item(13, 7).Value = 0.08417903999999998
$ws.Cells.Item(13, 8).Value = 1.14681178
$ws.Cells.Item(13, 13).Value = 0.03254535661374008
$ws.Cells.Item(13, 14).Value = 0.6719166098999868

